$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = 0.027509316802024841
$ws.Range("C13").Value = 0.038904048502445221
$ws.Range("B14").Value = 0.12835006415843964
$ws.Range("C14").Value = 0.18151441216468811
$ws.Range("B15").Value = 0.29355588555336
$ws.Range("C15").Value = 0.41515073180198669
$ws.Range("B16").Value = 0.48323655128479004
$ws.Range("C16").Value = 0.6833997368812561
$ws.Range("B17").Value = 0.97035729885101318
$ws.Range("C17").Value = 0.46429526805877686
$ws.Range("B18").Value = 1.62506902217865
$ws.Range("C18").Value = 0.16536343097686768
$ws.Range("B19").Value = 2.1426148414611816
$ws.Range("C19").Value = 0.81972277164459229
$ws.Range("B20").Value = 2.2663326263427734
$ws.Range("C20").Value = 1.110451340675354
$ws.Range("B21").Value = 1.932188868522644
$ws.Range("C21").Value = 0.86573046445846558
$ws.Range("B22").Value = 1.3124734163284302
$ws.Range("C22").Value = 0.22443512082099915
$ws.Range("B23").Value = 0.74116766452789307
$ws.Range("C23").Value = 0.42936956882476807
$ws.Range("B24").Value = 0.50957530736923218
$ws.Range("C24").Value = 0.72064834833145142
$ws.Range("B25").Value = 0.52736860513687134
$ws.Range("C25").Value = 0.74581176042556763
$ws.Range("B26").Value = 0.54398411512374878
$ws.Range("C26").Value = 0.76930969953536987
$ws.Range("B27").Value = 0.52024203538894653
$ws.Range("C27").Value = 0.73573338985443115
$ws.Range("B28").Value = 0.42965078353881836
$ws.Range("C28").Value = 0.607617974281311
$ws.Range("B29").Value = 0.27816802263259888
$ws.Range("C29").Value = 0.3933890163898468
$ws.Range("B30").Value = 0.11285752803087234
$ws.Range("C30").Value = 0.15960465371608734
$ws.Range("B31").Value = 0.005979481153190136
$ws.Range("C31").Value = 0.00845626275986433
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 0
$ws.Range("B39").Value = 0
$ws.Range("C39").Value = 0
$ws.Range("B40").Value = 0
$ws.Range("C40").Value = 0
$ws.Range("B41").Value = 0
$ws.Range("C41").Value = 0
$ws.Range("B42").Value = 0
$ws.Range("C42").Value = 0
$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("B44").Value = 0
$ws.Range("C44").Value = 0
$ws.Range("B45").Value = 0.36665940284729004
$ws.Range("C45").Value = 0.51853471994400024
$ws.Range("B46").Value = 1.6046019792556763
$ws.Range("C46").Value = 2.2692496776580811
$ws.Range("B47").Value = 3.2143280506134033
$ws.Range("C47").Value = 4.5457463264465332
$ws.Range("B48").Value = 4.4125714302062988
$ws.Range("C48").Value = 6.2403182983398438
$ws.Range("B49").Value = 5.3093714714050293
$ws.Range("C49").Value = 5.59274959564209
$ws.Range("B50").Value = 5.8198337554931641
$ws.Range("C50").Value = 2.4632678031921387
$ws.Range("B51").Value = 5.8727316856384277
$ws.Range("C51").Value = 1.9025601148605347
$ws.Range("B52").Value = 5.5966920852661133
$ws.Range("C52").Value = 5.6848325729370117
$ws.Range("B53").Value = 5.2206921577453613
$ws.Range("C53").Value = 7.383173942565918
$ws.Range("B54").Value = 4.7322049140930176
$ws.Range("C54").Value = 6.6923484802246094
$ws.Range("B55").Value = 3.5441751480102539
$ws.Range("C55").Value = 5.0122208595275879
$ws.Range("B56").Value = 2.081629753112793
$ws.Range("C56").Value = 2.94386887550354
$ws.Range("B57").Value = 0.82804673910140991
$ws.Range("C57").Value = 1.1710349321365356
$ws.Range("B58").Value = 0.12790371477603912
$ws.Range("C58").Value = 0.18088318407535553
$ws.Range("B59").Value = 0
$ws.Range("C59").Value = 0
$ws.Range("B60").Value = 0
$ws.Range("C60").Value = 0
$ws.Range("B61").Value = 0.615452766418457
$ws.Range("C61").Value = 0.87038165330886841
$ws.Range("B62").Value = 2.5719561576843262
$ws.Range("C62").Value = 1.7403228282928467
$ws.Range("B63").Value = 5.8038134574890137
$ws.Range("C63").Value = 1.7852746248245239
$ws.Range("B64").Value = 9.2452163696289063
$ws.Range("C64").Value = 1.0302048921585083
$ws.Range("B65").Value = 11.554998397827148
$ws.Range("C65").Value = 0.082054644823074341
$ws.Range("B66").Value = 11.75188159942627
$ws.Range("C66").Value = 1.0337026119232178
$ws.Range("B67").Value = 9.65302562713623
$ws.Range("C67").Value = 1.4676423072814941
$ws.Range("B68").Value = 6.01956033706665
$ws.Range("C68").Value = 1.307201623916626
$ws.Range("B69").Value = 2.3323037624359131
$ws.Range("C69").Value = 0.75929349660873413
$ws.Range("B70").Value = 0.14618940651416779
$ws.Range("C70").Value = 0.20674304664134979
$ws.Range("B71").Value = 0
$ws.Range("C71").Value = 0
$ws.Range("B72").Value = 0
$ws.Range("C72").Value = 0
